$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from AC1 into the new header cells, then set header labels
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record (Wins/Losses/Ties) for every data row
for ($row = 2; $row -le 43; $row++) {
    $ws.Range("AD$row").Value = 67
    $ws.Range("AE$row").Value = 94
    $ws.Range("AF$row").Value = 0
}
